# Integrationstests.xlsx - Kleine Typos in Integrationstests gefixt, Specs für Surface eingetragen

$wb = $excel.ActiveWorkbook

# --- Sheet "Tests": Typo-Fixes -------------------------------------------
$wsTests = $wb.Worksheets.Item("Tests")

# Zeile 38 (B38): "Play"/"Pause" waren vertauscht.
$wsTests.Range("B38").Value = 'Das "Pause"-Symbol wird durch ein "Play"-Symbol ersetzt.'

# Zeile 119 (A119): "Play"/"Pause" waren vertauscht.
$wsTests.Range("A119").Value = 'Der Nutzer klickt auf das "Pause"-Symbol, bevor das Spiel beendet ist.'

# --- Neues Blatt "Gerätespecs": Testgerät-Daten eintragen -----------------
$wsSpecs = $wb.Worksheets.Item("Gerätespecs")

$wsSpecs.Range("B1").Value = "Leonid Surface"

$wsSpecs.Range("A2").Value = "OS"
$wsSpecs.Range("B2").Value = "Windows 11 Pro"

$wsSpecs.Range("A3").Value = "Browser"
$wsSpecs.Range("B3").Value = "Microsoft Edge Version 120  (64 Bit)"

$wsSpecs.Range("A4").Value = "Anmerkungen"
$wsSpecs.Range("B4").Value = "mit Touchdisplay"

# --- Selektionen wie im finalen Dokument ----------------------------------
$wsSpecs.Select() | Out-Null
$wsSpecs.Range("B4").Select() | Out-Null

$wsTests.Select() | Out-Null
$wsTests.Range("B120").Select() | Out-Null
